$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to DAMSLTag (column I) and DialogAct (column J) following SGNN re-run
$updates = @(
    @{ Row = 40; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 42; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 51; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 55; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 77; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 104; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 118; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 122; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 141; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 145; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 147; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 153; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 155; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 157; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 158; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 159; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 160; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 161; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 163; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 168; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 169; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 179; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 195; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 196; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 197; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 198; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 199; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 203; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 211; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 212; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 213; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 214; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 231; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 235; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 243; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 258; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 267; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 281; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 284; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 288; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 302; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

